# make interval function di sisi mahasiswa saat menampilan data absensi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TANGGAL column (C2:C10) so every row shows the same interval date
$ws.Range("C2").Value = 43748
$ws.Range("C3").Value = 43748
$ws.Range("C4").Value = 43748
$ws.Range("C5").Value = 43748
$ws.Range("C6").Value = 43748
$ws.Range("C7").Value = 43748
$ws.Range("C8").Value = 43748
$ws.Range("C9").Value = 43748
$ws.Range("C10").Value = 43748

# Update the active selection to C4
$ws.Range("C4").Select()
